# Update the "Ich spiel..." sentences to "Ich spiele..." throughout column C,
# set row 12's height (as if text wrapped to a taller row), and move the
# active selection down to the second block of rows (A14:E25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 25; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value()
    if ($val -eq "Ich spiel mit meinem Hund") {
        $cell.Value = "Ich spiele mit meinem Hund"
    } elseif ($val -eq "Ich spiel mit meiner Katze") {
        $cell.Value = "Ich spiele mit meiner Katze"
    }
}

$ws.Rows.Item(12).RowHeight = 30.75

$ws.Range("A14:E25").Select() | Out-Null
